$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2023-11-16 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-17 Friday", 2) | Out-Null

# Update the division problems in the table
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "70÷6="
$t.Cell(1,2).Range.Text = "92÷7="
$t.Cell(1,3).Range.Text = "49÷6="
$t.Cell(1,4).Range.Text = "10÷5="
$t.Cell(1,5).Range.Text = "93÷3="

$t.Cell(5,1).Range.Text = "51÷6="
$t.Cell(5,2).Range.Text = "69÷2="
$t.Cell(5,3).Range.Text = "20÷9="
$t.Cell(5,4).Range.Text = "18÷5="
$t.Cell(5,5).Range.Text = "78÷3="

$t.Cell(9,1).Range.Text = "69÷5="
$t.Cell(9,2).Range.Text = "96÷8="
$t.Cell(9,3).Range.Text = "78÷3="
$t.Cell(9,4).Range.Text = "69÷5="
$t.Cell(9,5).Range.Text = "15÷9="

$t.Cell(13,1).Range.Text = "14÷8="
$t.Cell(13,2).Range.Text = "38÷3="
$t.Cell(13,3).Range.Text = "73÷7="
$t.Cell(13,4).Range.Text = "50÷5="
$t.Cell(13,5).Range.Text = "87÷2="

$t.Cell(17,1).Range.Text = "84÷9="
$t.Cell(17,2).Range.Text = "12÷3="
$t.Cell(17,3).Range.Text = "68÷7="
$t.Cell(17,4).Range.Text = "67÷7="
$t.Cell(17,5).Range.Text = "19÷6="
